$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 58247
$ws.Range("B2").Value = "Alícia Almeida"
$ws.Range("C2").Value = "Operacoes"
$ws.Range("D2").Value = "Viagem de negocios"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45085
$ws.Range("G2").Value = 2654.08

# Row 3
$ws.Range("A3").Value = 62492
$ws.Range("B3").Value = "João Guilherme Ramos"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 45086
$ws.Range("G3").Value = 5288.07

# Row 4
$ws.Range("A4").Value = 94840
$ws.Range("B4").Value = "Pietra Moraes"
$ws.Range("C4").Value = "Vendas"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45102
$ws.Range("G4").Value = 7883.34

# Row 5
$ws.Range("A5").Value = 29978
$ws.Range("B5").Value = "Dr. Anthony Marques"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45083
$ws.Range("G5").Value = 5849.76

# Row 6
$ws.Range("A6").Value = 11165
$ws.Range("B6").Value = "Helena Aparecida"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45087
$ws.Range("G6").Value = 2441.78

# Row 7
$ws.Range("A7").Value = 36234
$ws.Range("B7").Value = "Davi Luiz Sousa"
$ws.Range("D7").Value = "Consulta medica"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45086
$ws.Range("G7").Value = 4565.23

# Row 8
$ws.Range("A8").Value = 89671
$ws.Range("B8").Value = "Sophie Aragão"
$ws.Range("C8").Value = "Juridico"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 45086
$ws.Range("G8").Value = 9249.32

# Row 9
$ws.Range("A9").Value = 78341
$ws.Range("B9").Value = "Ana Cecília Martins"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 45098
$ws.Range("G9").Value = 3355.31

# Row 10
$ws.Range("A10").Value = 17758
$ws.Range("B10").Value = "Pedro Souza"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Viagem de negocios"
$ws.Range("F10").Value = 45100
$ws.Range("G10").Value = 2139.19

# Row 11
$ws.Range("A11").Value = 68785
$ws.Range("B11").Value = "Vinícius Mendonça"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45104
$ws.Range("G11").Value = 9841.82
